$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 115117.74
$ws.Range("I15").Value = 115117.74
$ws.Range("K15").Value = 345353.22
$ws.Range("M15").Value = -345184.22
$ws.Range("H74").Value = 9530678
$ws.Range("J74").Value = 10239.6
$ws.Range("L74").Value = 10239.6
$ws.Range("N74").Value = -12111.6
$ws.Range("H77").Value = 9530678
$ws.Range("J77").Value = 10239.6
$ws.Range("L77").Value = 51198
$ws.Range("N77").Value = -60558
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H115").Value = 903.55554
$ws.Range("I115").Value = 733.1429000000001
$ws.Range("K115").Value = 2199.4287
$ws.Range("M115").Value = -632.4287000000004
$ws.Range("H118").Value = 634.5714
$ws.Range("I118").Value = 259.5
$ws.Range("J118").Value = 1134.6666
$ws.Range("K118").Value = 778.5
$ws.Range("L118").Value = 3403.9998
$ws.Range("M118").Value = 878.5
$ws.Range("N118").Value = -6717.9998
$ws.Range("H128").Value = 24363.637
$ws.Range("I128").Value = 67000
$ws.Range("J128").Value = 22333.334
$ws.Range("K128").Value = 67000
$ws.Range("L128").Value = 22333.334
$ws.Range("M128").Value = -62020
$ws.Range("N128").Value = -32293.334
$ws.Range("H132").Value = 271538.84
$ws.Range("I132").Value = 324737.28
$ws.Range("K132").Value = 974211.8400000001
$ws.Range("M132").Value = -971681.8400000001
$ws.Range("H138").Value = 3042.8042
$ws.Range("J138").Value = 4029.6985
$ws.Range("L138").Value = 12089.0955
$ws.Range("N138").Value = -22369.0955

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4105.6274
$ws.Range("I2").Value = 2092.1025
$ws.Range("J2").Value = 10649.583
$ws.Range("K2").Value = 2092.1025
$ws.Range("L2").Value = 10649.583
$ws.Range("M2").Value = -1979.1025
$ws.Range("N2").Value = -10875.583
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H45").Value = 3269.0667
$ws.Range("I45").Value = 3202.5715
$ws.Range("J45").Value = 4200
$ws.Range("K45").Value = 3202.5715
$ws.Range("L45").Value = 4200
$ws.Range("M45").Value = -2825.5715
$ws.Range("N45").Value = -4954
$ws.Range("H54").Value = 36888.668
$ws.Range("I54").Value = 33333
$ws.Range("J54").Value = 44000
$ws.Range("K54").Value = 33333
$ws.Range("L54").Value = 44000
$ws.Range("M54").Value = -32564
$ws.Range("N54").Value = -45538
$ws.Range("H116").Value = 4105.6274
$ws.Range("I116").Value = 2092.1025
$ws.Range("J116").Value = 10649.583
$ws.Range("K116").Value = 2092.1025
$ws.Range("L116").Value = 10649.583
$ws.Range("M116").Value = 201.8975
$ws.Range("N116").Value = -15237.583
$ws.Range("H122").Value = 4147.636
$ws.Range("I122").Value = 3212.7896
$ws.Range("K122").Value = 9638.3688
$ws.Range("M122").Value = -7188.3688
$ws.Range("H132").Value = 1530821.4
$ws.Range("I132").Value = 2087743.1
$ws.Range("K132").Value = 6263229.300000001
$ws.Range("M132").Value = -6260699.300000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4105.6274
$ws.Range("I3").Value = 2092.1025
$ws.Range("J3").Value = 10649.583
$ws.Range("K3").Value = 2092.1025
$ws.Range("L3").Value = 10649.583
$ws.Range("M3").Value = -1978.1025
$ws.Range("N3").Value = -10877.583
$ws.Range("H99").Value = 7342.4053
$ws.Range("I99").Value = 7030.1914
$ws.Range("J99").Value = 7800.9688
$ws.Range("K99").Value = 7030.1914
$ws.Range("L99").Value = 7800.9688
$ws.Range("M99").Value = -5532.1914
$ws.Range("N99").Value = -10796.9688
$ws.Range("H105").Value = 3628.9
$ws.Range("I105").Value = 3628.9
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3628.9
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1881.9
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 1020441.8
$ws.Range("I134").Value = 1374141.6
$ws.Range("K134").Value = 4122424.8
$ws.Range("M134").Value = -4119889.8

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12823803
$ws.Range("I16").Value = 15153132
$ws.Range("K16").Value = 15153132
$ws.Range("M16").Value = -15152845
$ws.Range("H58").Value = 17930
$ws.Range("I58").Value = 17113
$ws.Range("J58").Value = 18474.666
$ws.Range("K58").Value = 17113
$ws.Range("L58").Value = 18474.666
$ws.Range("M58").Value = -16910
$ws.Range("N58").Value = -18880.666
$ws.Range("H113").Value = 12823803
$ws.Range("I113").Value = 15153132
$ws.Range("K113").Value = 15153132
$ws.Range("M113").Value = -15150962
$ws.Range("H134").Value = 10752.7
$ws.Range("J134").Value = 15345.7
$ws.Range("L134").Value = 46037.10000000001
$ws.Range("N134").Value = -51107.10000000001
$ws.Range("H136").Value = 17930
$ws.Range("I136").Value = 17113
$ws.Range("J136").Value = 18474.666
$ws.Range("K136").Value = 51339
$ws.Range("L136").Value = 55423.99800000001
$ws.Range("M136").Value = -48789
$ws.Range("N136").Value = -60523.99800000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 543.7778
$ws.Range("I8").Value = 543.7778
$ws.Range("K8").Value = 1631.3334
$ws.Range("M8").Value = -1492.3334
$ws.Range("H17").Value = 3818.625
$ws.Range("I17").Value = 6816.3335
$ws.Range("J17").Value = 2020
$ws.Range("K17").Value = 20449.0005
$ws.Range("L17").Value = 6060
$ws.Range("M17").Value = -20280.0005
$ws.Range("N17").Value = -6398

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 211.875
$ws.Range("I2").Value = 259.7
$ws.Range("J2").Value = 132.16667
$ws.Range("K2").Value = 259.7
$ws.Range("L2").Value = 132.16667
$ws.Range("M2").Value = -146.7
$ws.Range("N2").Value = -358.16667
$ws.Range("H80").Value = 9023.916999999999
$ws.Range("J80").Value = 10842.375
$ws.Range("L80").Value = 10842.375
$ws.Range("N80").Value = -12838.375
$ws.Range("H83").Value = 9023.916999999999
$ws.Range("J83").Value = 10842.375
$ws.Range("L83").Value = 54211.875
$ws.Range("N83").Value = -64195.875
$ws.Range("H132").Value = 5905.136
$ws.Range("J132").Value = 6342.75
$ws.Range("L132").Value = 19028.25
$ws.Range("N132").Value = -24088.25
$ws.Range("H136").Value = 23678.473
$ws.Range("J136").Value = 23678.473
$ws.Range("L136").Value = 71035.41900000001
$ws.Range("N136").Value = -76135.41900000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3475.0645
$ws.Range("I40").Value = 3024.5833
$ws.Range("J40").Value = 5019.5713
$ws.Range("K40").Value = 3024.5833
$ws.Range("L40").Value = 5019.5713
$ws.Range("M40").Value = -2888.5833
$ws.Range("N40").Value = -5291.5713
$ws.Range("H132").Value = 5022.6665
$ws.Range("I132").Value = 3801.6511
$ws.Range("K132").Value = 11404.9533
$ws.Range("M132").Value = -8874.953300000001
$ws.Range("H136").Value = 9300.954
$ws.Range("I136").Value = 7228.615
$ws.Range("J136").Value = 12294.333
$ws.Range("K136").Value = 21685.845
$ws.Range("L136").Value = 36882.999
$ws.Range("M136").Value = -19135.845
$ws.Range("N136").Value = -41982.999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1846.6818
$ws.Range("I100").Value = 1569
$ws.Range("J100").Value = 2790.8
$ws.Range("K100").Value = 3138
$ws.Range("L100").Value = 5581.6
$ws.Range("M100").Value = -2597
$ws.Range("N100").Value = -6663.6
$ws.Range("H126").Value = 2800.2
$ws.Range("I126").Value = 1557.2727
$ws.Range("K126").Value = 4671.8181
$ws.Range("M126").Value = -2201.8181
$ws.Range("H132").Value = 9143.781000000001
$ws.Range("I132").Value = 9001.85
$ws.Range("J132").Value = 9380.333000000001
$ws.Range("K132").Value = 27005.55
$ws.Range("L132").Value = 28140.999
$ws.Range("M132").Value = -24475.55
$ws.Range("N132").Value = -33200.999
$ws.Range("H136").Value = 8629624
$ws.Range("I136").Value = 14294142
$ws.Range("J136").Value = 9704.434999999999
$ws.Range("K136").Value = 42882426
$ws.Range("L136").Value = 29113.305
$ws.Range("M136").Value = -42879876
$ws.Range("N136").Value = -34213.305
